$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "tc149" - only the saved selection moves (no data change)
# ---------------------------------------------------------------------------
$wsTc149 = $wb.Worksheets.Item("tc149")
$wsTc149.Activate()
$wsTc149.Range("G1").Select()

# ---------------------------------------------------------------------------
# Sheet "tcasst011" (ASST011 test case) - rework the test data: one row
# (Asset011) becomes three rows (Asset01/02/03) and the "Type" column header
# is highlighted.
# ---------------------------------------------------------------------------
$wsAsst011 = $wb.Worksheets.Item("tcasst011")
$wsAsst011.Activate()

# Highlight the "Type" header cell (C1) the same way the other completed
# test-case sheets do.
$wsTc149.Range("I1").Copy()
$wsAsst011.Range("C1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# Clone row 2's formatting down into rows 3 & 4 before filling in values.
$wsAsst011.Range("A2:K2").Copy()
$wsAsst011.Range("A3:K4").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# Row 2: Asset01
$wsAsst011.Range("A2").Value = "Asset01"
$wsAsst011.Range("B2").Value = "1"
$wsAsst011.Range("C2").Value = "HeatBath"
$wsAsst011.Range("D2").Value = "1-Manufacturer"
$wsAsst011.Range("E2").Value = "1-location"
$wsAsst011.Range("F2").Value = "Model-11"
$wsAsst011.Range("G2").Value = "5"
$wsAsst011.Range("H2").Value = "cu ft"
$wsAsst011.Range("I2").Value = "3"
$wsAsst011.Range("J2").Value = "Years"
$wsAsst011.Range("K2").Value = "ASST011-Test"

# Row 3: Asset02
$wsAsst011.Range("A3").Value = "Asset02"
$wsAsst011.Range("B3").Value = "2"
$wsAsst011.Range("C3").Value = "HeatBath"
$wsAsst011.Range("D3").Value = "2-Manufacturer"
$wsAsst011.Range("E3").Value = "2-location"
$wsAsst011.Range("F3").Value = "Model-12"
$wsAsst011.Range("G3").Value = "5"
$wsAsst011.Range("H3").Value = "cu ft"
$wsAsst011.Range("I3").Value = "3"
$wsAsst011.Range("J3").Value = "Years"
$wsAsst011.Range("K3").Value = "ASST011-Test"

# Row 4: Asset03
$wsAsst011.Range("A4").Value = "Asset03"
$wsAsst011.Range("B4").Value = "3"
$wsAsst011.Range("C4").Value = "Sterilizer"
$wsAsst011.Range("D4").Value = "3-Manufacturer"
$wsAsst011.Range("E4").Value = "3-location"
$wsAsst011.Range("F4").Value = "Model-13"
$wsAsst011.Range("G4").Value = "5"
$wsAsst011.Range("H4").Value = "cu ft"
$wsAsst011.Range("I4").Value = "3"
$wsAsst011.Range("J4").Value = "Years"
$wsAsst011.Range("K4").Value = "ASST011-Test"

$wsAsst011.Range("J13").Select()

# ---------------------------------------------------------------------------
# Sheet "tcasst014" (ASST014 test case) - rework the test data: one row
# (Asset019) becomes three rows (Asset01/02/03) and the "Manufacturer"
# column header is highlighted.
# ---------------------------------------------------------------------------
$wsAsst014 = $wb.Worksheets.Item("tcasst014")
$wsAsst014.Activate()

# Highlight the "Manufacturer" header cell (D1).
$wsTc149.Range("I1").Copy()
$wsAsst014.Range("D1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# Clone row 2's formatting down into rows 3 & 4 before filling in values.
$wsAsst014.Range("A2:K2").Copy()
$wsAsst014.Range("A3:K4").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# Row 2: Asset01
$wsAsst014.Range("A2").Value = "Asset01"
$wsAsst014.Range("B2").Value = "1"
$wsAsst014.Range("C2").Value = "HeatBath"
$wsAsst014.Range("D2").Value = "1-Manufacturer"
$wsAsst014.Range("E2").Value = "1-location"
$wsAsst014.Range("F2").Value = "Model-14"
$wsAsst014.Range("G2").Value = "5"
$wsAsst014.Range("H2").Value = "cu ft"
$wsAsst014.Range("I2").Value = "3"
$wsAsst014.Range("J2").Value = "Years"
$wsAsst014.Range("K2").Value = "ASST014- Test"

# Row 3: Asset02
$wsAsst014.Range("A3").Value = "Asset02"
$wsAsst014.Range("B3").Value = "2"
$wsAsst014.Range("C3").Value = "HeatBath"
$wsAsst014.Range("D3").Value = "1-Man"
$wsAsst014.Range("E3").Value = "1-locate"
$wsAsst014.Range("F3").Value = "Model-15"
$wsAsst014.Range("G3").Value = "5"
$wsAsst014.Range("H3").Value = "cu ft"
$wsAsst014.Range("I3").Value = "3"
$wsAsst014.Range("J3").Value = "Years"
$wsAsst014.Range("K3").Value = "ASST014- Test"

# Row 4: Asset03
$wsAsst014.Range("A4").Value = "Asset03"
$wsAsst014.Range("B4").Value = "3"
$wsAsst014.Range("C4").Value = "HeatBath"
$wsAsst014.Range("D4").Value = "1-Manuf"
$wsAsst014.Range("E4").Value = "1-locatio"
$wsAsst014.Range("F4").Value = "Model-16"
$wsAsst014.Range("G4").Value = "5"
$wsAsst014.Range("H4").Value = "cu ft"
$wsAsst014.Range("I4").Value = "3"
$wsAsst014.Range("J4").Value = "Years"
$wsAsst014.Range("K4").Value = "ASST014- Test"

$wsAsst014.Range("D11").Select()
